$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (Förändrad) for rows 2 through 205 changes from 45175 to 45177
$ws.Range("C2:C205").Value = 45177
